# Rotate the data rows 15-18 upward by one (row16->15, row17->16, row18->17),
# wrapping the original row 15 content around into row 18.
# (The four species-observation records occupying rows 15-18 were
# re-sorted; every field of each row moves together as a unit.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All columns that carry data anywhere in rows 15-18.
$cols = @("A","B","D","E","F","G","H","I","P","Q","R","S","T","U","V","W", `
          "Y","Z","AA","AB","AC","AD","AE","AF","AG","AT","AW","AX","AY")

# Columns whose content is a plain date-looking string (e.g. "2026-01-31")
# that Excel would otherwise silently reinterpret as a date serial number
# when it is written back through COM.
$dateTextCols = @("Y","AA")

function Get-RowValues($ws, [int]$row, $cols) {
    $vals = @{}
    foreach ($col in $cols) {
        $vals[$col] = $ws.Range("$col$row").Value2
    }
    return $vals
}

function Set-RowValues($ws, [int]$row, $cols, $oldVals, $newVals, $dateTextCols) {
    foreach ($col in $cols) {
        $old = $oldVals[$col]
        $new = $newVals[$col]

        if (($old -eq $null) -and ($new -eq $null)) {
            continue
        }
        if (($old -ne $null) -and ($new -ne $null) -and ("$old" -ceq "$new")) {
            continue
        }

        $rng = $ws.Range("$col$row")

        if ($new -eq $null) {
            # Target wants a blank cell; clear it out entirely.
            $rng.ClearContents()
        } elseif (($dateTextCols -contains $col)) {
            # Force the destination cell to stay plain text so the
            # "yyyy-mm-dd" string isn't converted into a date value.
            $rng.NumberFormat = "@"
            $rng.Value2 = $new
            $rng.ClearFormats()
        } else {
            $rng.Value2 = $new
            if ("$new" -eq "") {
                # Keep the cell present (touched) even though it is blank,
                # matching source rows that had an explicit-but-empty cell.
                $rng.Font.Bold = $false
            }
        }
    }
}

# Capture the original contents of the four rows before overwriting anything.
$row15 = Get-RowValues $ws 15 $cols
$row16 = Get-RowValues $ws 16 $cols
$row17 = Get-RowValues $ws 17 $cols
$row18 = Get-RowValues $ws 18 $cols

# Write them back shifted: 16->15, 17->16, 18->17, 15->18.
Set-RowValues $ws 15 $cols $row15 $row16 $dateTextCols
Set-RowValues $ws 16 $cols $row16 $row17 $dateTextCols
Set-RowValues $ws 17 $cols $row17 $row18 $dateTextCols
Set-RowValues $ws 18 $cols $row18 $row15 $dateTextCols

Write-Output "Rows 15-18 rotated."
